$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.343.71"
$ws.Range("E2").Value = "  +3.12%  "

$ws.Range("D3").Value = "3.123.07"
$ws.Range("E3").Value = "  +1.05%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'219.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.43%  "

$ws.Range("D6").Value = "'624.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.39%  "

$ws.Range("D7").Value = "'0.993"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +23.85%  "

$ws.Range("D8").Value = "'0.380"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.43%  "

$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").Value = "3.118.23"
$ws.Range("E10").Value = "  +1.03%  "

$ws.Range("D11").Value = "'0.725"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +21.97%  "

$ws.Range("D12").Value = "'0.192"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.57%  "

$ws.Range("E13").Value = "  +6.31%  "

$ws.Range("D14").Value = "'34.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.00%  "

$ws.Range("D15").Value = "'5.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.69%  "

$ws.Range("D16").Value = "91.037.62"
$ws.Range("E16").Value = "  +3.54%  "

$ws.Range("D17").Value = "3.691.48"
$ws.Range("E17").Value = "  +1.00%  "

$ws.Range("D18").Value = "3.118.98"
$ws.Range("E18").Value = "  +1.44%  "

$ws.Range("D19").Value = "'3.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +12.16%  "

$ws.Range("D20").Value = "'0.0000220"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.56%  "

$ws.Range("D21").Value = "'14.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.17%  "

$ws.Range("D22").Value = "'436.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.48%  "

$ws.Range("D23").Value = "'8.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.34%  "

$ws.Range("E24").Value = "  +5.07%  "

$ws.Range("D25").Value = "'6.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.89%  "

$ws.Range("D26").Value = "'12.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.69%  "

$ws.Range("D27").Value = "'86.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.28%  "

$ws.Range("D28").Value = "3.279.04"
$ws.Range("E28").Value = "  +1.46%  "

$ws.Range("E29").Value = "  -0.39%  "

$ws.Range("E30").Value = "  -1.05%  "

$ws.Range("D31").Value = "'9.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +13.61%  "

$ws.Range("D32").Value = "'529.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.77%  "

$ws.Range("D33").Value = "'0.896"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -17.18%  "

$ws.Range("D34").Value = "'3.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.65%  "

$ws.Range("D35").Value = "'7.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.03%  "

$ws.Range("E36").Value = "  +12.55%  "

$ws.Range("D37").Value = "'1.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.53%  "

$ws.Range("D38").Value = "'23.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.52%  "

$ws.Range("D39").Value = "'1.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.82%  "

$ws.Range("D40").Value = "'0.0902"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +30.11%  "

$ws.Range("E41").Value = "  +0.22%  "

$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("D43").Value = "'0.154"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +16.80%  "

$ws.Range("D44").Value = "'0.401"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.34%  "

$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.87%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").Value = "'149.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.29%  "

$ws.Range("D48").Value = "'44.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.68%  "

$ws.Range("D49").Value = "'1.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.60%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'168.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.61%  "

$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").Value = "'4.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.13%  "
